$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row cells (AP1:AS1), bold/centered style to match existing headers
$ws.Range("AP1").Value = "ACCEPTED"
$ws.Range("AQ1").Value = "PAID"
$ws.Range("AR1").Value = "HOLD"
$ws.Range("AS1").Value = "REJECTED"

# Copy header style (bold, border, centered) from an existing header cell
$ws.Range("AO1").Copy()
$ws.Range("AP1:AS1").PasteSpecial(-4122)  # xlPasteFormats

# Row 2
$ws.Range("AI2").Value = "HOLD"
$ws.Range("AJ2").Value = "HOLD"
$ws.Range("AP2").Value = $false
$ws.Range("AQ2").Value = $false
$ws.Range("AR2").Value = $true
$ws.Range("AS2").Value = $false

# Row 3
$ws.Range("AI3").Value = "HOLD"
$ws.Range("AJ3").Value = "HOLD"
$ws.Range("AP3").Value = $false
$ws.Range("AQ3").Value = $false
$ws.Range("AR3").Value = $true
$ws.Range("AS3").Value = $false

# Rows 4-9
for ($r = 4; $r -le 9; $r++) {
    $ws.Range("AP$r").Value = $false
    $ws.Range("AQ$r").Value = $false
    $ws.Range("AR$r").Value = $false
    $ws.Range("AS$r").Value = $false
}
